$wb = $excel.ActiveWorkbook

# The workbook has two sheets that carry duplicate exhibition data:
# "展览" (sheet 1) and "全部类型" (sheet 4). Both need the same F-column
# (ticket/attendance count) updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F9").Value = 1577
    $ws.Range("F12").Value = 62
    $ws.Range("F13").Value = 486
}
